$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "lesson #40 lesson video record on youtube"
# Lesson #39 (row 43) gets pushed out two days and its recording link is filled in;
# lesson #40 (row 44) gets its real name and recording link filled in too.

# Row 43 - lesson #39 ("Spring intro"): date 2021-02-16 -> 2021-02-18, add YouTube link.
$ws.Range("E43").Value = 44245
$ws.Rows.Item(43).RowHeight = 13.85

$ws.Hyperlinks.Add($ws.Range("F43"), "https://youtu.be/E3LZMBqVjQQ", "", "", "https://youtu.be/E3LZMBqVjQQ")
$ws.Range("F43").Value = "https://youtu.be/E3LZMBqVjQQ "
# Hyperlinks.Add stamps the built-in blue/underlined Hyperlink font; put the
# sheet's normal (non-hyperlink) look back, matching the other link cells.
$ws.Range("F43").Font.Name = "Cambria"
$ws.Range("F43").Font.Size = 11
$ws.Range("F43").Font.Underline = 0
$ws.Range("F43").Font.ColorIndex = 1

# Row 44 - lesson #40: placeholder name replaced with "Spring MVC", add YouTube link.
$ws.Range("C44").Value = "Spring MVC"

$ws.Hyperlinks.Add($ws.Range("F44"), "https://youtu.be/AJnCbBv3o3o", "", "", "https://youtu.be/AJnCbBv3o3o")
$ws.Range("F44").Value = "https://youtu.be/AJnCbBv3o3o "
$ws.Range("F44").Font.Name = "Cambria"
$ws.Range("F44").Font.Size = 11
$ws.Range("F44").Font.Underline = 0
$ws.Range("F44").Font.ColorIndex = 1

# Author left the selection on F47 when the edit was saved.
$ws.Range("F47").Select() | Out-Null
